$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8943026715494398
$ws.Range("C2").Value = 0.2207033664204516
$ws.Range("D2").Value = 0.6463507405400435
$ws.Range("E2").Value = 0.2640199051264958
$ws.Range("G2").Value = 0.5641648832301627
$ws.Range("H2").Value = 0.6794157543653654
$ws.Range("I2").Value = 0.4555822578606126
$ws.Range("J2").Value = 0.1377449535323407
$ws.Range("N2").Value = 0.9928180756836795
$ws.Range("O2").Value = 2.462085459092208

# Row 3
$ws.Range("B3").Value = 0.7983501396883526
$ws.Range("C3").Value = 0.1938668191324666
$ws.Range("D3").Value = 0.6355131684779849
$ws.Range("E3").Value = 0.2587149794968155
$ws.Range("G3").Value = 0.5601569542302798
$ws.Range("H3").Value = 0.6819553233815725
$ws.Range("I3").Value = 0.4607534365761055
$ws.Range("J3").Value = 0.1342381605419618
$ws.Range("N3").Value = 0.9955897625118837
$ws.Range("O3").Value = 2.458583285427864

# Row 4
$ws.Range("B4").Value = 0.7394098925317394
$ws.Range("C4").Value = 0.177351592814091
$ws.Range("D4").Value = 0.629192979444241
$ws.Range("E4").Value = 0.2556003478742213
$ws.Range("G4").Value = 0.5581432135588784
$ws.Range("H4").Value = 0.6838758657685275
$ws.Range("I4").Value = 0.464247973045449
$ws.Range("J4").Value = 0.1321636387888461
$ws.Range("N4").Value = 0.9976915983850958
$ws.Range("O4").Value = 2.458125525397378

# Row 5
$ws.Range("B5").Value = 0.7153865102837926
$ws.Range("C5").Value = 0.1706124497135306
$ws.Range("D5").Value = 0.6267015522837198
$ws.Range("E5").Value = 0.2543669864242872
$ws.Range("G5").Value = 0.5574347905124597
$ws.Range("H5").Value = 0.6847493231125554
$ws.Range("I5").Value = 0.4657522070288493
$ws.Range("J5").Value = 0.1313380249694447
$ws.Range("N5").Value = 0.9986489032584842
$ws.Range("O5").Value = 2.458364039722198

# Row 6
$ws.Range("B6").Value = 0.7113972041525471
$ws.Range("C6").Value = 0.169492884188827
$ws.Range("D6").Value = 0.6262929346762576
$ws.Range("E6").Value = 0.2541643549792099
$ws.Range("G6").Value = 0.5573239263876246
$ws.Range("H6").Value = 0.6848998449600572
$ws.Range("I6").Value = 0.4660068227286089
$ws.Range("J6").Value = 0.1312021262411349
$ws.Range("N6").Value = 0.9988139565576901
$ws.Range("O6").Value = 2.458429300790897

# Row 7
$ws.Range("B7").Value = 0.739085922007348
$ws.Range("C7").Value = 0.1772607425436377
$ws.Range("D7").Value = 0.6291590385523591
$ws.Range("E7").Value = 0.2555835690486532
$ws.Range("G7").Value = 0.5581332055804182
$ws.Range("H7").Value = 0.6838872778116496
$ws.Range("I7").Value = 0.4642679351993486
$ws.Range("J7").Value = 0.1321524242345475
$ws.Range("N7").Value = 0.9977041005666436
$ws.Range("O7").Value = 2.458127021800522

# Row 8
$ws.Range("B8").Value = 0.8612243182632255
$ws.Range("C8").Value = 0.211458116961154
$ws.Range("D8").Value = 0.6425446429332737
$ws.Range("E8").Value = 0.2621611736721121
$ws.Range("G8").Value = 0.5626899939842218
$ws.Range("H8").Value = 0.6802164126988259
$ws.Range("I8").Value = 0.4572989042172786
$ws.Range("J8").Value = 0.136519477329557
$ws.Range("N8").Value = 0.9936908509086138
$ws.Range("O8").Value = 2.460526266988722

# Row 9
$ws.Range("B9").Value = 1.100483224427819
$ws.Range("C9").Value = 0.278208651519833
$ws.Range("D9").Value = 0.6714434191805481
$ws.Range("E9").Value = 0.2761917838373265
$ws.Range("G9").Value = 0.5751861823226534
$ws.Range("H9").Value = 0.6758853839070298
$ws.Range("I9").Value = 0.4461732482453264
$ws.Range("J9").Value = 0.1457084857370745
$ws.Range("N9").Value = 0.9889866788372785
$ws.Range("O9").Value = 2.478690546316784

# Row 10
$ws.Range("B10").Value = 1.276054193981395
$ws.Range("C10").Value = 0.327048228606543
$ws.Range("D10").Value = 0.6942918716734141
$ws.Range("E10").Value = 0.2871921387892371
$ws.Range("G10").Value = 0.5865568451016401
$ws.Range("H10").Value = 0.6744541169232434
$ws.Range("I10").Value = 0.4395566809312434
$ws.Range("J10").Value = 0.1528431953021112
$ws.Range("N10").Value = 0.9874508872261316
$ws.Range("O10").Value = 2.500287798479349

# Row 11
$ws.Range("B11").Value = 1.355868692195656
$ws.Range("C11").Value = 0.3492203268741605
$ws.Range("D11").Value = 0.7050377088406208
$ws.Range("E11").Value = 0.2923472839407637
$ws.Range("G11").Value = 0.5922093832505482
$ws.Range("H11").Value = 0.6741837869796541
$ws.Range("I11").Value = 0.43688668906346
$ws.Range("J11").Value = 0.1561728113328513
$ws.Range("N11").Value = 0.9871673522067397
$ws.Range("O11").Value = 2.511915234135813

# Row 12
$ws.Range("B12").Value = 1.386083335704257
$ws.Range("C12").Value = 0.3576095034179048
$ws.Range("D12").Value = 0.7091574596218493
$ws.Range("E12").Value = 0.2943211318280774
$ws.Range("G12").Value = 0.5944191735630824
$ws.Range("H12").Value = 0.6741362124868431
$ws.Range("I12").Value = 0.43592466728715
$ws.Range("J12").Value = 0.1574457560473945
$ws.Range("N12").Value = 0.9871195165037534
$ws.Range("O12").Value = 2.516578205261425

# Row 13
$ws.Range("B13").Value = 1.379576515015515
$ws.Range("C13").Value = 0.3558030581858702
$ws.Range("D13").Value = 0.7082679512135712
$ws.Range("E13").Value = 0.2938950631574144
$ws.Range("G13").Value = 0.5939401692659061
$ws.Range("H13").Value = 0.6741440208930527
$ws.Range("I13").Value = 0.4361296718968148
$ws.Range("J13").Value = 0.1571710666662796
$ws.Range("N13").Value = 0.9871271735054137
$ws.Range("O13").Value = 2.51556238041664

# Row 14
$ws.Range("B14").Value = 1.358354666505193
$ws.Range("C14").Value = 0.3499106504040697
$ws.Range("D14").Value = 0.7053756308467882
$ws.Range("E14").Value = 0.292509238633194
$ws.Range("G14").Value = 0.5923897936278451
$ws.Range("H14").Value = 0.6741787746176726
$ws.Range("I14").Value = 0.4368065590410168
$ws.Range("J14").Value = 0.1562772948614963
$ws.Range("N14").Value = 0.9871622244968705
$ws.Range("O14").Value = 2.512293646617309

# Row 15
$ws.Range("B15").Value = 1.345354401327143
$ws.Range("C15").Value = 0.3463004678143875
$ws.Range("D15").Value = 0.7036105800573296
$ws.Range("E15").Value = 0.2916632075592958
$ws.Range("G15").Value = 0.5914491765141889
$ws.Range("H15").Value = 0.6742071992343028
$ws.Range("I15").Value = 0.4372275641628072
$ws.Range("J15").Value = 0.1557314088394861
$ws.Range("N15").Value = 0.9871914425625903
$ws.Range("O15").Value = 2.51032532143455

# Row 16
$ws.Range("B16").Value = 1.270836907493788
$ws.Range("C16").Value = 0.3255982841482705
$ws.Range("D16").Value = 0.6935966817945314
$ws.Range("E16").Value = 0.2868582760341738
$ws.Range("G16").Value = 0.5861971203920291
$ws.Range("H16").Value = 0.6744794489657266
$ws.Range("I16").Value = 0.4397380254223222
$ws.Range("J16").Value = 0.1526272884166815
$ws.Range("N16").Value = 0.9874777561780377
$ws.Range("O16").Value = 2.499564252230243

# Row 17
$ws.Range("B17").Value = 1.22510792011775
$ws.Range("C17").Value = 0.3128862809202246
$ws.Range("D17").Value = 0.6875435711278612
$ws.Range("E17").Value = 0.2839492800822327
$ws.Range("G17").Value = 0.5830982832127489
$ws.Range("H17").Value = 0.6747440161050804
$ws.Range("I17").Value = 0.4413652916297011
$ws.Range("J17").Value = 0.1507445343819285
$ws.Range("N17").Value = 0.9877596038349026
$ws.Range("O17").Value = 2.493424894754469

# Row 18
$ws.Range("B18").Value = 1.198800868629746
$ws.Range("C18").Value = 0.3055704360302798
$ws.Range("D18").Value = 0.6840951148107592
$ws.Range("E18").Value = 0.2822903187554431
$ws.Range("G18").Value = 0.5813610741728041
$ws.Range("H18").Value = 0.6749320246509001
$ws.Range("I18").Value = 0.4423332368329049
$ws.Range("J18").Value = 0.1496695298948083
$ws.Range("N18").Value = 0.9879607958535814
$ws.Range("O18").Value = 2.490063340136061

# Row 19
$ws.Range("B19").Value = 1.189892954096877
$ws.Range("C19").Value = 0.3030927013616065
$ws.Range("D19").Value = 0.6829332196726341
$ws.Range("E19").Value = 0.2817310649489144
$ws.Range("G19").Value = 0.5807806328528216
$ws.Range("H19").Value = 0.6750018348413391
$ws.Range("I19").Value = 0.4426664547835095
$ws.Range("J19").Value = 0.1493069093097432
$ws.Range("N19").Value = 0.9880356330896234
$ws.Range("O19").Value = 2.488954290790105

# Row 20
$ws.Range("B20").Value = 1.229976373003694
$ws.Range("C20").Value = 0.3142399368091446
$ws.Range("D20").Value = 0.6881845065210825
$ws.Range("E20").Value = 0.2842574764269017
$ws.Range("G20").Value = 0.5834234837965226
$ws.Range("H20").Value = 0.6747121432401286
$ws.Range("I20").Value = 0.4411887548583628
$ws.Range("J20").Value = 0.1509441384511803
$ws.Range("N20").Value = 0.9877255571298917
$ws.Range("O20").Value = 2.494060877637907

# Row 21
$ws.Range("B21").Value = 1.364588305509415
$ws.Range("C21").Value = 0.3516415843284904
$ws.Range("D21").Value = 0.706223804538439
$ws.Range("E21").Value = 0.2929156998771916
$ws.Range("G21").Value = 0.5928432936260606
$ws.Range("H21").Value = 0.6741670791979715
$ws.Range("I21").Value = 0.4366064084527927
$ws.Range("J21").Value = 0.1565394890526903
$ws.Range("N21").Value = 0.987150314702987
$ws.Range("O21").Value = 2.513246693556539

# Row 22
$ws.Range("B22").Value = 1.452509572627207
$ws.Range("C22").Value = 0.37604522446901
$ws.Range("D22").Value = 0.7183080526554591
$ws.Range("E22").Value = 0.2987008650378513
$ws.Range("G22").Value = 0.5994037271857025
$ws.Range("H22").Value = 0.6741302464699714
$ws.Range("I22").Value = 0.4338975361572039
$ws.Range("J22").Value = 0.1602668575981738
$ws.Range("N22").Value = 0.9871212938558358
$ws.Range("O22").Value = 2.527301024647784

# Row 23
$ws.Range("B23").Value = 1.405589956733252
$ws.Range("C23").Value = 0.3630243819451948
$ws.Range("D23").Value = 0.7118315381880507
$ws.Range("E23").Value = 0.2956016411008022
$ws.Range("G23").Value = 0.595865238623503
$ws.Range("H23").Value = 0.6741206663003823
$ws.Range("I23").Value = 0.4353170923144951
$ws.Range("J23").Value = 0.1582710376022334
$ws.Range("N23").Value = 0.9871050893897717
$ws.Range("O23").Value = 2.519661089222637

# Row 24
$ws.Range("B24").Value = 1.227775397804294
$ws.Range("C24").Value = 0.3136279724546966
$ws.Range("D24").Value = 0.687894641327091
$ws.Range("E24").Value = 0.2841180989359557
$ws.Range("G24").Value = 0.5832763224736226
$ws.Range("H24").Value = 0.6747264411196454
$ws.Range("I24").Value = 0.4412684661941455
$ws.Range("J24").Value = 0.1508538743540271
$ws.Range("N24").Value = 0.9877408276681905
$ws.Range("O24").Value = 2.493772826374226

# Row 25
$ws.Range("B25").Value = 1.035790532808505
$ws.Range("C25").Value = 0.2601854231375569
$ws.Range("D25").Value = 0.663341745513975
$ws.Range("E25").Value = 0.2722747426092269
$ws.Range("G25").Value = 0.5714224831459092
$ws.Range("H25").Value = 0.6767497846496866
$ws.Range("I25").Value = 0.4489101523926351
$ws.Range("J25").Value = 0.1431554333788583
$ws.Range("N25").Value = 0.98992137528397
$ws.Range("O25").Value = 2.472330875022124
